# Fruta / hortaliza, semanal
#
# This adds one new weekly price record. It is modeled as inserting a new
# row at row 52 (which pushes the existing rows 52-129 down to 53-130,
# preserving all of their data/formatting), then filling the newly blank
# row 52 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 52; rows 52..129 shift down to 53..130.
$ws.Rows("52").Insert()

# Populate the new row 52 with the new weekly record.
$ws.Range("A52").Value = 3
$ws.Range("B52").Value = "Femacal de La Calera"
$ws.Range("C52").Value = "Coquimbo"
$ws.Range("D52").Value = 44495
$ws.Range("E52").Value = 5
$ws.Range("F52").Value = "Fruta"
$ws.Range("G52").Value = 100101
$ws.Range("H52").Value = "Berries"
$ws.Range("I52").Value = 100101001
$ws.Range("J52").Value = "Arándano (blue)"
$ws.Range("K52").Value = "Sin especificar"
$ws.Range("L52").Value = "Primera"
$ws.Range("M52").Value = 48
$ws.Range("N52").Value = 10000
$ws.Range("O52").Value = 10000
$ws.Range("P52").Value = 10000
$ws.Range("Q52").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R52").Value = "Provincia de Quillota"
$ws.Range("S52").Value = 6667
$ws.Range("T52").Value = 1.5
